$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22; existing rows 22-102 shift down to 23-103
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with data
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44687
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100102
$ws.Range("H22").Value = "Cítricos"
$ws.Range("I22").Value = 100102004
$ws.Range("J22").Value = "Mandarina"
$ws.Range("K22").Value = "Murcott"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 13000
$ws.Range("P22").Value = 12500
$ws.Range("Q22").Value = "$/bandeja 18 kilos"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 694
$ws.Range("T22").Value = 18
